$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Teacher"
$ws.Range("D1").Value = "Nurse"
$ws.Range("F1").Value = 7
$ws.Range("G1").Value = "3,1,5,2,6,4"
$ws.Range("E2").Value = "Teacher"
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = "1,6,2,3"
$ws.Range("C3").Value = "Teacher"
$ws.Range("E3").Value = "Engineer"
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = "5,4,6,3"
$ws.Range("C4").Value = "Doctor"
$ws.Range("D4").Value = "Teacher"
$ws.Range("E4").Value = "Engineer"
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = "5,3,1,4"
$ws.Range("C5").Value = "Engineer"
$ws.Range("D5").Value = "Doctor"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = "1,2,4,5"
$ws.Range("D6").Value = "Nurse"
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 5
$ws.Range("D7").Value = "Nurse"
$ws.Range("F7").Value = 13
$ws.Range("G7").Value = "6,4,5,1"
$ws.Range("C8").Value = "Teacher"
$ws.Range("D8").Value = "Nurse"
$ws.Range("E8").Value = "Doctor"
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = "4,2,6,3,5"
$ws.Range("E9").Value = "Teacher"
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = "2,1,6,3"
$ws.Range("C10").Value = "Doctor"
$ws.Range("D10").Value = "Engineer"
$ws.Range("E10").Value = "Nurse"
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = "2,5,6,3"
$ws.Range("C11").Value = "Doctor"
$ws.Range("D11").Value = "Teacher"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = "2,3,5,4,6,1"
$ws.Range("D12").Value = "Nurse"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 14
$ws.Range("G12").Value = "5,4,3,1,6"
$ws.Range("C13").Value = "Doctor"
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = "5,3,1,6"
$ws.Range("C14").Value = "Teacher"
$ws.Range("E14").Value = "Doctor"
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = "5,2,3,1,4"
$ws.Range("D15").Value = "Engineer"
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = "3,2,4,6"
$ws.Range("C16").Value = "Teacher"
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = "2,1,4,5,6,3"
$ws.Range("C17").Value = "Engineer"
$ws.Range("D17").Value = "Doctor"
$ws.Range("E17").Value = "Teacher"
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = "5,2,3,1,6"
$ws.Range("C18").Value = "Doctor"
$ws.Range("D18").Value = "Nurse"
$ws.Range("E18").Value = "Teacher"
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = "5,2,6,4,1"
$ws.Range("C19").Value = "Nurse"
$ws.Range("E19").Value = "Teacher"
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = "4,5,2,1,6"
$ws.Range("C20").Value = "Engineer"
$ws.Range("D20").Value = "Nurse"
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = "2,1,4,5,3,6"
